$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update row 2 values
$ws.Range("D2").Value = 1115970.9
$ws.Range("E2").Value = 2000
$ws.Range("F2").Value = 30000
$ws.Range("H2").Value = 11272651347.38185
$ws.Range("I2").Value = 18536.07700107334
$ws.Range("J2").Value = -11853029973.01902

# Update row 3 values
$ws.Range("D3").Value = 1115970.9
$ws.Range("E3").Value = 2000
$ws.Range("F3").Value = 30000
$ws.Range("H3").Value = 11272146584.22761
$ws.Range("I3").Value = 18645.40180227152
$ws.Range("J3").Value = -12441723946.87802

# Remove rows 4 and 5 entirely (data no longer needed)
$ws.Range("A4:J5").Delete()
